# Daily attendance processing - 2025-10-07 12:35:49
#
# Normalize the "Recorded By" (column G) entries: when a row's recorder
# list begins with the literal token "System," (i.e. "System" was logged
# first in the comma-separated list of recorders), move it to the end of
# the list instead, preserving the original casing/order of every other
# token. Rows whose recorder list does not start with "System," are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System,")) {
        $parts = $val -split ",\s*"
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
